# Applies updated cryptocurrency price/volume data to sheet1
# (GitHub Actions scheduled refresh of cryptos list)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces text interpretation so decimal-looking
# price strings (e.g. '1.006') are not auto-converted to numbers.
$ws.Range("D2").Value = "'28.437.16"
$ws.Range("E2").Value = "  -3.51%  "
$ws.Range("D3").Value = "'1.957.58"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  -0.90%  "
$ws.Range("D5").Value = "'320.54"
$ws.Range("E5").Value = "  -2.76%  "
$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D7").Value = "'0.4763"
$ws.Range("E7").Value = "  -5.28%  "
$ws.Range("D8").Value = "'0.4055"
$ws.Range("E8").Value = "  -4.18%  "
$ws.Range("D9").Value = "'53.52"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").Value = "'0.08425"
$ws.Range("E10").Value = "  -5.69%  "
$ws.Range("D11").Value = "'1.058"
$ws.Range("E11").Value = "  -4.78%  "
$ws.Range("D12").Value = "'22.43"
$ws.Range("E12").Value = "  -3.27%  "
$ws.Range("D13").Value = "'1.962.38"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("E14").Value = "  -4.53%  "
$ws.Range("E15").Value = "  -4.88%  "
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "'90.15"
$ws.Range("E17").Value = "  -4.43%  "
$ws.Range("D18").Value = "'0.00001068"
$ws.Range("E18").Value = "  -3.85%  "
$ws.Range("D19").Value = "'0.06590"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("D20").Value = "'18.52"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").Value = "'5.817"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "'28.442.85"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24").Value = "'11.52"
$ws.Range("E24").Value = "  -4.52%  "
$ws.Range("D25").Value = "'2.291"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").Value = "'2.162.62"
$ws.Range("E26").Value = "  -3.75%  "
$ws.Range("D27").Value = "'155.57"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").Value = "'20.19"
$ws.Range("E28").Value = "  -2.79%  "
$ws.Range("D29").Value = "'5.916"
$ws.Range("E29").Value = "  -5.86%  "
$ws.Range("D30").Value = "'2.154"
$ws.Range("E30").Value = "  -6.38%  "
$ws.Range("D31").Value = "'123.65"
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("E32").Value = "  -7.59%  "
$ws.Range("D33").Value = "'0.09597"
$ws.Range("E33").Value = "  -3.48%  "
$ws.Range("D34").Value = "'1.451"
$ws.Range("E34").Value = "  -5.91%  "
$ws.Range("D35").Value = "'5.603"
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("D36").Value = "'3.661"
$ws.Range("E36").Value = "  -3.54%  "
$ws.Range("D37").Value = "'8.964"
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("D38").Value = "'0.02327"
$ws.Range("E38").Value = "  -5.31%  "
$ws.Range("D39").Value = "'0.06218"
$ws.Range("E39").Value = "  -2.67%  "
$ws.Range("D40").Value = "'1.240"
$ws.Range("E40").Value = "  -4.27%  "
$ws.Range("D41").Value = "'0.6209"
$ws.Range("E41").Value = "  -4.94%  "
$ws.Range("D42").Value = "'11.13"
$ws.Range("E42").Value = "  -4.32%  "
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").Value = "'0.1924"
$ws.Range("E44").Value = "  -5.66%  "
$ws.Range("D45").Value = "'1.359"
$ws.Range("E45").Value = "  +3.70%  "
$ws.Range("D46").Value = "'0.5955"
$ws.Range("E46").Value = "  -5.80%  "
$ws.Range("D47").Value = "'12.97"
$ws.Range("E47").Value = "  -3.67%  "
$ws.Range("D48").Value = "'2.061"
$ws.Range("E48").Value = "  -6.68%  "
$ws.Range("D49").Value = "'3.388"
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("D50").Value = "'0.00000000328"
$ws.Range("E50").Value = "  -4.44%  "
$ws.Range("D51").Value = "'0.06819"
$ws.Range("E51").Value = "  -1.93%  "
